$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (row 1) ---
$ws.Range("B1").Value = "time (min)"
$ws.Range("C1").Value = "part"
$ws.Range("D1").Value = "comment"
$ws.Range("G1").Value = "time (hour)"

# --- New rows 3-5 ---
$ws.Range("A3").Value = 210824
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "read more part 1 and do ex 1.1 and 1.2"

$ws.Range("A4").Value = 210831
$ws.Range("B4").Value = 80
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "ex 1.3 and 1.4 and respective material"

$ws.Range("A5").Value = 210908
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "ex 1.5 and 1.6 and respective material"

# --- Row 2: update comment text (time/part/date stay the same) ---
$ws.Range("D2").Value = "read material of parts 0 and 1 and exercise guide and install stuff"

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 10.416666666666666
$ws.Columns.Item(4).ColumnWidth = 56.583333333333336
$ws.Columns.Item(7).ColumnWidth = 14.916666666666666

# --- Selection ---
[void]$ws.Range("J13").Select()
